$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 2258.4546
$ws.Range("I9").Value = 2056.8333
$ws.Range("J9").Value = 2500.4
$ws.Range("K9").Value = 2056.8333
$ws.Range("L9").Value = 2500.4
$ws.Range("M9").Value = -1887.8333
$ws.Range("N9").Value = -2838.4
$ws.Range("H15").Value = 1527.18
$ws.Range("I15").Value = 1527.18
$ws.Range("K15").Value = 4581.54
$ws.Range("M15").Value = -4412.54
$ws.Range("H28").Value = 2452.5833
$ws.Range("I28").Value = 1765.909
$ws.Range("K28").Value = 1765.909
$ws.Range("M28").Value = -1280.909
$ws.Range("H33").Value = 429
$ws.Range("I33").Value = 443.22223
$ws.Range("K33").Value = 443.22223
$ws.Range("M33").Value = -214.22223
$ws.Range("H39").Value = 737.2222
$ws.Range("J39").Value = 999
$ws.Range("L39").Value = 2997
$ws.Range("N39").Value = -3589
$ws.Range("H43").Value = 5115.7646
$ws.Range("I43").Value = 1999
$ws.Range("J43").Value = 5783.643
$ws.Range("K43").Value = 1999
$ws.Range("L43").Value = 5783.643
$ws.Range("M43").Value = -1930
$ws.Range("N43").Value = -5921.643
$ws.Range("H64").Value = 6665.5
$ws.Range("I64").Value = 4297.5
$ws.Range("J64").Value = 7849.5
$ws.Range("K64").Value = 4297.5
$ws.Range("L64").Value = 7849.5
$ws.Range("M64").Value = -4049.5
$ws.Range("N64").Value = -8345.5
$ws.Range("H67").Value = 6665.5
$ws.Range("I67").Value = 4297.5
$ws.Range("J67").Value = 7849.5
$ws.Range("K67").Value = 4297.5
$ws.Range("L67").Value = 7849.5
$ws.Range("M67").Value = -3439.5
$ws.Range("N67").Value = -9565.5
$ws.Range("H70").Value = 4375372
$ws.Range("J70").Value = 5561474
$ws.Range("L70").Value = 16684422
$ws.Range("N70").Value = -16684962
$ws.Range("H73").Value = 4375372
$ws.Range("J73").Value = 5561474
$ws.Range("L73").Value = 16684422
$ws.Range("N73").Value = -16686294
$ws.Range("H80").Value = 2253.2273
$ws.Range("I80").Value = 1193.8
$ws.Range("K80").Value = 3581.4
$ws.Range("M80").Value = -2583.4
$ws.Range("H83").Value = 2253.2273
$ws.Range("I83").Value = 1193.8
$ws.Range("K83").Value = 10744.2
$ws.Range("M83").Value = -5752.199999999999
$ws.Range("H88").Value = 236069.23
$ws.Range("I88").Value = 1500450
$ws.Range("J88").Value = 6181.8184
$ws.Range("K88").Value = 1500450
$ws.Range("L88").Value = 6181.8184
$ws.Range("M88").Value = -1500044
$ws.Range("N88").Value = -6993.8184
$ws.Range("H91").Value = 236069.23
$ws.Range("I91").Value = 1500450
$ws.Range("J91").Value = 6181.8184
$ws.Range("K91").Value = 1500450
$ws.Range("L91").Value = 6181.8184
$ws.Range("M91").Value = -1499046
$ws.Range("N91").Value = -8989.8184
$ws.Range("H100").Value = 6485.857
$ws.Range("I100").Value = 4248.75
$ws.Range("J100").Value = 9468.666999999999
$ws.Range("K100").Value = 4248.75
$ws.Range("L100").Value = 9468.666999999999
$ws.Range("M100").Value = -3707.75
$ws.Range("N100").Value = -10550.667
$ws.Range("H129").Value = 2203
$ws.Range("I129").Value = 1247.8
$ws.Range("K129").Value = 3743.4
$ws.Range("M129").Value = 1256.6
$ws.Range("H132").Value = 1684.6857
$ws.Range("I132").Value = 1606.5151
$ws.Range("K132").Value = 4819.5453
$ws.Range("M132").Value = -2289.5453
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -80119
$ws.Range("H135").Value = 3096.3
$ws.Range("I135").Value = 2804.1428
$ws.Range("J135").Value = 3778
$ws.Range("K135").Value = 25237.2852
$ws.Range("L135").Value = 34002
$ws.Range("M135").Value = -22702.2852
$ws.Range("N135").Value = -39072
$ws.Range("H137").Value = 3103.8772
$ws.Range("I137").Value = 2450.7273
$ws.Range("J137").Value = 3260.0652
$ws.Range("K137").Value = 7352.1819
$ws.Range("L137").Value = 9780.195599999999
$ws.Range("M137").Value = -4802.1819
$ws.Range("N137").Value = -14880.1956
$ws.Range("H138").Value = 2851.5476
$ws.Range("I138").Value = 1208.9231
$ws.Range("J138").Value = 3587.8965
$ws.Range("K138").Value = 3626.7693
$ws.Range("L138").Value = 10763.6895
$ws.Range("M138").Value = 1513.2307
$ws.Range("N138").Value = -21043.6895
$ws.Range("H141").Value = 2237.3076
$ws.Range("I141").Value = 1943.6364
$ws.Range("J141").Value = 3852.5
$ws.Range("K141").Value = 5830.9092
$ws.Range("L141").Value = 11557.5
$ws.Range("M141").Value = -650.9092000000001
$ws.Range("N141").Value = -21917.5

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 3274.2075
$ws.Range("I32").Value = 3067.673
$ws.Range("K32").Value = 3067.673
$ws.Range("M32").Value = -2780.673
$ws.Range("H45").Value = 1650.8235
$ws.Range("I45").Value = 1643.125
$ws.Range("J45").Value = 1774
$ws.Range("K45").Value = 1643.125
$ws.Range("L45").Value = 1774
$ws.Range("M45").Value = -1266.125
$ws.Range("N45").Value = -2528
$ws.Range("H61").Value = 9432.950000000001
$ws.Range("I61").Value = 4043.1333
$ws.Range("K61").Value = 4043.1333
$ws.Range("M61").Value = -3831.1333
$ws.Range("H74").Value = 22226622
$ws.Range("I74").Value = 37039416
$ws.Range("K74").Value = 37039416
$ws.Range("M74").Value = -37038542
$ws.Range("H77").Value = 22226622
$ws.Range("I77").Value = 37039416
$ws.Range("K77").Value = 185197080
$ws.Range("M77").Value = -185192712
$ws.Range("H80").Value = 79723.75
$ws.Range("J80").Value = 139998
$ws.Range("L80").Value = 139998
$ws.Range("N80").Value = -141994
$ws.Range("H83").Value = 79723.75
$ws.Range("J83").Value = 139998
$ws.Range("L83").Value = 419994
$ws.Range("N83").Value = -429978
$ws.Range("H97").Value = 1154.84
$ws.Range("I97").Value = 1323.375
$ws.Range("K97").Value = 1323.375
$ws.Range("M97").Value = -827.375
$ws.Range("H110").Value = 4455.6665
$ws.Range("I110").Value = 3605.5625
$ws.Range("K110").Value = 3605.5625
$ws.Range("M110").Value = -1560.5625
$ws.Range("H122").Value = 2293.2222
$ws.Range("I122").Value = 2142.375
$ws.Range("K122").Value = 6427.125
$ws.Range("M122").Value = -3977.125
$ws.Range("H132").Value = 11944.723
$ws.Range("I132").Value = 10785.286
$ws.Range("J132").Value = 16002.75
$ws.Range("K132").Value = 32355.858
$ws.Range("L132").Value = 48008.25
$ws.Range("M132").Value = -29825.858
$ws.Range("N132").Value = -53068.25
$ws.Range("H135").Value = 54765.43
$ws.Range("J135").Value = 54765.43
$ws.Range("L135").Value = 54765.43
$ws.Range("N135").Value = -64905.43
$ws.Range("H136").Value = 9432.950000000001
$ws.Range("I136").Value = 4043.1333
$ws.Range("K136").Value = 12129.3999
$ws.Range("M136").Value = -9579.3999
$ws.Range("H137").Value = 69998
$ws.Range("J137").Value = 69998
$ws.Range("L137").Value = 69998
$ws.Range("N137").Value = -80198

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 4876
$ws.Range("I20").Value = 1052
$ws.Range("K20").Value = 1052
$ws.Range("M20").Value = -805
$ws.Range("H82").Value = 6394
$ws.Range("I82").Value = 6394
$ws.Range("K82").Value = 6394
$ws.Range("M82").Value = -6011
$ws.Range("H85").Value = 6394
$ws.Range("I85").Value = 6394
$ws.Range("K85").Value = 6394
$ws.Range("M85").Value = -5068
$ws.Range("H97").Value = 17766.889
$ws.Range("I97").Value = 4985
$ws.Range("J97").Value = 33744.25
$ws.Range("K97").Value = 4985
$ws.Range("L97").Value = 33744.25
$ws.Range("M97").Value = -3994
$ws.Range("N97").Value = -35726.25
$ws.Range("H105").Value = 5658.0435
$ws.Range("I105").Value = 2424.3572
$ws.Range("K105").Value = 2424.3572
$ws.Range("M105").Value = -677.3571999999999
$ws.Range("H107").Value = 1784.125
$ws.Range("I107").Value = 1457.3846
$ws.Range("K107").Value = 1457.3846
$ws.Range("M107").Value = 462.6153999999999
$ws.Range("H111").Value = 52000
$ws.Range("J111").Value = 52000
$ws.Range("L111").Value = 52000
$ws.Range("N111").Value = -60180
$ws.Range("H122").Value = 80199.8
$ws.Range("I122").Value = 91999
$ws.Range("J122").Value = 77250
$ws.Range("K122").Value = 91999
$ws.Range("L122").Value = 77250
$ws.Range("M122").Value = -87099
$ws.Range("N122").Value = -87050
$ws.Range("H134").Value = 2441.3333
$ws.Range("I134").Value = 2429.8
$ws.Range("K134").Value = 7289.400000000001
$ws.Range("M134").Value = -4754.400000000001
$ws.Range("H135").Value = 46713.832
$ws.Range("J135").Value = 46713.832
$ws.Range("L135").Value = 46713.832
$ws.Range("N135").Value = -56853.832

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1815.6666
$ws.Range("I22").Value = 423.16666
$ws.Range("J22").Value = 3208.1667
$ws.Range("K22").Value = 423.16666
$ws.Range("L22").Value = 3208.1667
$ws.Range("M22").Value = -73.16665999999998
$ws.Range("N22").Value = -3908.1667
$ws.Range("H31").Value = 30988.61
$ws.Range("I31").Value = 3724.1667
$ws.Range("J31").Value = 69479.586
$ws.Range("K31").Value = 3724.1667
$ws.Range("L31").Value = 69479.586
$ws.Range("M31").Value = -3429.1667
$ws.Range("N31").Value = -70069.586
$ws.Range("H34").Value = 30988.61
$ws.Range("I34").Value = 3724.1667
$ws.Range("J34").Value = 69479.586
$ws.Range("K34").Value = 3724.1667
$ws.Range("L34").Value = 69479.586
$ws.Range("M34").Value = -3522.1667
$ws.Range("N34").Value = -69883.586
$ws.Range("H62").Value = 9257.200000000001
$ws.Range("I62").Value = 4093.3333
$ws.Range("J62").Value = 17003
$ws.Range("K62").Value = 4093.3333
$ws.Range("L62").Value = 17003
$ws.Range("M62").Value = -3469.3333
$ws.Range("N62").Value = -18251
$ws.Range("H65").Value = 9257.200000000001
$ws.Range("I65").Value = 4093.3333
$ws.Range("J65").Value = 17003
$ws.Range("K65").Value = 20466.6665
$ws.Range("L65").Value = 85015
$ws.Range("M65").Value = -17346.6665
$ws.Range("N65").Value = -91255
$ws.Range("H99").Value = 2514.8462
$ws.Range("I99").Value = 1799
$ws.Range("K99").Value = 1799
$ws.Range("M99").Value = -301
$ws.Range("H100").Value = 69999
$ws.Range("J100").Value = 69999
$ws.Range("L100").Value = 69999
$ws.Range("N100").Value = -72163
$ws.Range("H105").Value = 2641.5715
$ws.Range("I105").Value = 1098.4
$ws.Range("K105").Value = 1098.4
$ws.Range("M105").Value = 648.5999999999999
$ws.Range("H122").Value = 9922.272000000001
$ws.Range("I122").Value = 5399.3335
$ws.Range("K122").Value = 16198.0005
$ws.Range("M122").Value = -13748.0005
$ws.Range("H126").Value = 2514.8462
$ws.Range("I126").Value = 1799
$ws.Range("K126").Value = 5397
$ws.Range("M126").Value = -2927
$ws.Range("H132").Value = 4207.5625
$ws.Range("I132").Value = 2787.25
$ws.Range("K132").Value = 8361.75
$ws.Range("M132").Value = -5831.75
$ws.Range("H133").Value = 62419.223
$ws.Range("J133").Value = 62419.223
$ws.Range("L133").Value = 62419.223
$ws.Range("N133").Value = -67479.223
$ws.Range("H138").Value = 62800.188
$ws.Range("J138").Value = 62800.188
$ws.Range("L138").Value = 62800.188
$ws.Range("N138").Value = -73080.18799999999

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H4").Value = 9026445
$ws.Range("I4").Value = 3913378.5
$ws.Range("J4").Value = 20786496
$ws.Range("K4").Value = 11740135.5
$ws.Range("L4").Value = 62359488
$ws.Range("M4").Value = -11740023.5
$ws.Range("N4").Value = -62359712
$ws.Range("H5").Value = 1757468.8
$ws.Range("I5").Value = 1950
$ws.Range("K5").Value = 5850
$ws.Range("M5").Value = -5738
$ws.Range("H14").Value = 2291.4443
$ws.Range("I14").Value = 2291.4443
$ws.Range("K14").Value = 6874.3329
$ws.Range("M14").Value = -6701.3329
$ws.Range("H68").Value = 2184.8572
$ws.Range("J68").Value = 1948.5
$ws.Range("L68").Value = 5845.5
$ws.Range("N68").Value = -7467.5
$ws.Range("H71").Value = 2184.8572
$ws.Range("J71").Value = 1948.5
$ws.Range("L71").Value = 17536.5
$ws.Range("N71").Value = -25648.5
$ws.Range("H75").Value = 66668708
$ws.Range("J75").Value = 3631.5715
$ws.Range("L75").Value = 10894.7145
$ws.Range("N75").Value = -12890.7145
$ws.Range("H78").Value = 66668708
$ws.Range("J78").Value = 3631.5715
$ws.Range("L78").Value = 32684.1435
$ws.Range("N78").Value = -42668.1435
$ws.Range("H80").Value = 7984.8887
$ws.Range("I80").Value = 5992
$ws.Range("K80").Value = 17976
$ws.Range("M80").Value = -17040
$ws.Range("H83").Value = 7984.8887
$ws.Range("I83").Value = 5992
$ws.Range("K83").Value = 53928
$ws.Range("M83").Value = -49248
$ws.Range("H86").Value = 3489.6155
$ws.Range("I86").Value = 2817
$ws.Range("J86").Value = 3788.5557
$ws.Range("K86").Value = 8451
$ws.Range("L86").Value = 11365.6671
$ws.Range("M86").Value = -7265
$ws.Range("N86").Value = -13737.6671
$ws.Range("H88").Value = 13344
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 6000
$ws.Range("M88").Value = -5572
$ws.Range("H89").Value = 3489.6155
$ws.Range("I89").Value = 2817
$ws.Range("J89").Value = 3788.5557
$ws.Range("K89").Value = 25353
$ws.Range("L89").Value = 34097.0013
$ws.Range("M89").Value = -19425
$ws.Range("N89").Value = -45953.0013
$ws.Range("H91").Value = 13344
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 6000
$ws.Range("M91").Value = -4518
$ws.Range("H97").Value = 373.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 373.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1120.5
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -2112.5
$ws.Range("H98").Value = 1396.3636
$ws.Range("I98").Value = 850.8
$ws.Range("J98").Value = 1851
$ws.Range("K98").Value = 2552.4
$ws.Range("L98").Value = 5553
$ws.Range("M98").Value = -1054.4
$ws.Range("N98").Value = -8549
$ws.Range("H107").Value = 3127318
$ws.Range("I107").Value = 3351.25
$ws.Range("J107").Value = 5209962.5
$ws.Range("K107").Value = 10053.75
$ws.Range("L107").Value = 15629887.5
$ws.Range("M107").Value = -8133.75
$ws.Range("N107").Value = -15633727.5
$ws.Range("H116").Value = 2004916.5
$ws.Range("I116").Value = 2669878
$ws.Range("K116").Value = 8009634
$ws.Range("M116").Value = -8006192
$ws.Range("H126").Value = 4999
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 6656858.5
$ws.Range("I131").Value = 8930010
$ws.Range("J131").Value = 5720854.5
$ws.Range("K131").Value = 26790030
$ws.Range("L131").Value = 17162563.5
$ws.Range("M131").Value = -26784990
$ws.Range("N131").Value = -17172643.5
$ws.Range("H132").Value = 2806
$ws.Range("I132").Value = 1023
$ws.Range("J132").Value = 3499.389
$ws.Range("K132").Value = 9207
$ws.Range("L132").Value = 31494.501
$ws.Range("M132").Value = -6677
$ws.Range("N132").Value = -36554.501
$ws.Range("H135").Value = 1757468.8
$ws.Range("I135").Value = 1950
$ws.Range("K135").Value = 17550
$ws.Range("M135").Value = -15015
$ws.Range("H137").Value = 61276.06
$ws.Range("J137").Value = 102666.3
$ws.Range("L137").Value = 307998.9
$ws.Range("N137").Value = -318198.9

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H43").Value = 500
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H57").Value = 6571.4287
$ws.Range("I57").Value = 1000
$ws.Range("K57").Value = 1000
$ws.Range("M57").Value = -180
$ws.Range("H70").Value = 10097.272
$ws.Range("I70").Value = 10341.111
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 10341.111
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -10071.111
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 10097.272
$ws.Range("I73").Value = 10341.111
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 10341.111
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -9405.111000000001
$ws.Range("N73").Value = -10872
$ws.Range("H80").Value = 4163.6
$ws.Range("I80").Value = 4329.5
$ws.Range("K80").Value = 4329.5
$ws.Range("M80").Value = -3331.5
$ws.Range("H83").Value = 4163.6
$ws.Range("I83").Value = 4329.5
$ws.Range("K83").Value = 21647.5
$ws.Range("M83").Value = -16655.5
$ws.Range("H97").Value = 2560.4
$ws.Range("I97").Value = 3590
$ws.Range("J97").Value = 1874
$ws.Range("K97").Value = 3590
$ws.Range("L97").Value = 1874
$ws.Range("M97").Value = -3094
$ws.Range("N97").Value = -2866
$ws.Range("H107").Value = 1327.5625
$ws.Range("J107").Value = 2177.4285
$ws.Range("L107").Value = 2177.4285
$ws.Range("N107").Value = -6017.4285
$ws.Range("H132").Value = 5083.778
$ws.Range("I132").Value = 4276.8423
$ws.Range("K132").Value = 12830.5269
$ws.Range("M132").Value = -10300.5269

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 5541.5454
$ws.Range("I7").Value = 4012.7878
$ws.Range("J7").Value = 10127.818
$ws.Range("K7").Value = 4012.7878
$ws.Range("L7").Value = 10127.818
$ws.Range("M7").Value = -3900.7878
$ws.Range("N7").Value = -10351.818
$ws.Range("H22").Value = 4037.1538
$ws.Range("I22").Value = 923
$ws.Range("J22").Value = 9019.799999999999
$ws.Range("K22").Value = 923
$ws.Range("L22").Value = 9019.799999999999
$ws.Range("M22").Value = -628
$ws.Range("N22").Value = -9609.799999999999
$ws.Range("H24").Value = 9000
$ws.Range("J24").Value = 9000
$ws.Range("L24").Value = 9000
$ws.Range("N24").Value = -9686
$ws.Range("H27").Value = 4037.1538
$ws.Range("I27").Value = 923
$ws.Range("J27").Value = 9019.799999999999
$ws.Range("K27").Value = 923
$ws.Range("L27").Value = 9019.799999999999
$ws.Range("M27").Value = -816
$ws.Range("N27").Value = -9233.799999999999
$ws.Range("H34").Value = 11000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 11000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 11000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -11344
$ws.Range("H36").Value = 79571.664
$ws.Range("J36").Value = 79571.664
$ws.Range("L36").Value = 79571.664
$ws.Range("N36").Value = -80695.664
$ws.Range("H68").Value = 3154.2222
$ws.Range("I68").Value = 3236.5
$ws.Range("K68").Value = 3236.5
$ws.Range("M68").Value = -2487.5
$ws.Range("H71").Value = 3154.2222
$ws.Range("I71").Value = 3236.5
$ws.Range("K71").Value = 16182.5
$ws.Range("M71").Value = -12438.5
$ws.Range("H82").Value = 9246.041999999999
$ws.Range("I82").Value = 5255.6924
$ws.Range("K82").Value = 5255.6924
$ws.Range("M82").Value = -4894.6924
$ws.Range("H85").Value = 9246.041999999999
$ws.Range("I85").Value = 5255.6924
$ws.Range("K85").Value = 5255.6924
$ws.Range("M85").Value = -4007.6924
$ws.Range("H100").Value = 5084.0835
$ws.Range("J100").Value = 8762.4
$ws.Range("L100").Value = 8762.4
$ws.Range("N100").Value = -9844.4
$ws.Range("H126").Value = 5541.5454
$ws.Range("I126").Value = 4012.7878
$ws.Range("J126").Value = 10127.818
$ws.Range("K126").Value = 12038.3634
$ws.Range("L126").Value = 30383.454
$ws.Range("M126").Value = -9568.3634
$ws.Range("N126").Value = -35323.454
$ws.Range("H133").Value = 52492.75
$ws.Range("J133").Value = 53325
$ws.Range("L133").Value = 53325
$ws.Range("N133").Value = -58385
$ws.Range("H136").Value = 6769.6294
$ws.Range("I136").Value = 2714.611
$ws.Range("J136").Value = 14879.667
$ws.Range("K136").Value = 8143.833
$ws.Range("L136").Value = 44639.001
$ws.Range("M136").Value = -5593.833
$ws.Range("N136").Value = -49739.001
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 15301
$ws.Range("J62").Value = 15301
$ws.Range("L62").Value = 15301
$ws.Range("N62").Value = -16549
$ws.Range("H65").Value = 15301
$ws.Range("J65").Value = 15301
$ws.Range("L65").Value = 76505
$ws.Range("N65").Value = -82745
$ws.Range("H81").Value = 5631.8887
$ws.Range("I81").Value = 3969.5715
$ws.Range("J81").Value = 11450
$ws.Range("K81").Value = 7939.143
$ws.Range("L81").Value = 22900
$ws.Range("M81").Value = -6878.143
$ws.Range("N81").Value = -25022
$ws.Range("H84").Value = 5631.8887
$ws.Range("I84").Value = 3969.5715
$ws.Range("J84").Value = 11450
$ws.Range("K84").Value = 39695.715
$ws.Range("L84").Value = 114500
$ws.Range("M84").Value = -34391.715
$ws.Range("N84").Value = -125108
$ws.Range("H126").Value = 1738.4
$ws.Range("I126").Value = 1565.24
$ws.Range("J126").Value = 2171.3
$ws.Range("K126").Value = 4695.72
$ws.Range("L126").Value = 6513.900000000001
$ws.Range("M126").Value = -2225.72
$ws.Range("N126").Value = -11453.9
$ws.Range("H132").Value = 4953.125
$ws.Range("I132").Value = 2910.75
$ws.Range("J132").Value = 9037.875
$ws.Range("K132").Value = 8732.25
$ws.Range("L132").Value = 27113.625
$ws.Range("M132").Value = -6202.25
$ws.Range("N132").Value = -32173.625
$ws.Range("H139").Value = 69964.91
$ws.Range("I139").Value = 69487.5
$ws.Range("K139").Value = 69487.5
$ws.Range("M139").Value = -64347.5
